# Update the poster title's year from "(2024)" to "(2025)".
#
# The title placeholder ("Rectangle 2") holds a single paragraph whose first
# run reads "DSS5202 Sustainable Systems Analysis (2024)". We only touch the
# "(2024)" substring so the rest of the run (and the line break / author-name
# run / endParaRPr that follow it) are left completely untouched, matching
# how PowerPoint itself would split a run when a user edits a portion of it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$title = $s.Shapes.Title
$tr = $title.TextFrame.TextRange

$fullText = $tr.Text
$oldYear = "(2024)"
$newYear = "(2025)"

$startPos = $fullText.IndexOf($oldYear) + 1   # PowerPoint ranges are 1-based

if ($startPos -gt 0) {
    $yearRange = $tr.Characters($startPos, $oldYear.Length)
    $yearRange.Text = $newYear
}
